$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.434.56'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").Value = '1.937.59'
$ws.Range("E3").Value = '  -0.78%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '242.34'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("E6").Value = '  -1.03%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '56.35'
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  -3.47%  '
$ws.Range("E9").Value = '  -3.40%  '
$ws.Range("E10").Value = '  -3.72%  '
$ws.Range("E11").Value = '  -1.16%  '
$ws.Range("D12").Value = '2.222.38'
$ws.Range("E12").Value = '  -0.89%  '
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '21.36'
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  -2.27%  '
$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '0.799'
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  -3.93%  '
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '13.20'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  -2.47%  '
$ws.Range("E16").Value = '  -3.98%  '
$ws.Range("D17").Value = '1.935.81'
$ws.Range("E17").Value = '  -0.51%  '
$ws.Range("D18").Value = '36.398.45'
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("E19").Value = '  -1.74%  '
$ws.Range("E20").Value = '  -2.93%  '
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '226.31'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  -1.64%  '
$ws.Range("E22").Value = '  -2.81%  '
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("E24").Value = '  -5.48%  '
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("E26").Value = '  -4.53%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '159.30'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -2.84%  '
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '0.132'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  +8.47%  '
$ws.Range("E29").Value = '  -3.21%  '
$ws.Range("E30").Value = '  -0.80%  '
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '1.07'
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -6.69%  '
$ws.Range("E32").Value = '  -3.73%  '
$ws.Range("E33").Value = '  -4.11%  '
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '4.11'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  -4.64%  '
$ws.Range("E35").Value = '  -0.19%  '
$ws.Range("E36").Value = '  -0.34%  '
$ws.Range("E37").Value = '  -1.33%  '
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '2.19'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  +1.12%  '
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '3.19'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  +10.13%  '
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '0.0980'
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  -0.13%  '
$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '2.90'
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  +1.07%  '
$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '0.0208'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  -1.24%  '
$ws.Range("E43").Value = '  -4.05%  '
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '15.64'
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  -0.36%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '1.328.30'
$ws.Range("E45").Value = '  -1.26%  '
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '1.02'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  -2.37%  '
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '85.36'
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  -3.88%  '
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '7.05'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -4.17%  '
$ws.Range("E49").Value = '  -0.34%  '
$ws.Range("D50").Value = '2.114.70'
$ws.Range("E50").Value = '  -0.91%  '
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '3.46'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +14.49%  '
